# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "23.889.66"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.646.71"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "308.83"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3886"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.79%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3823"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "51.09"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.98%  "

$ws.Range("E10").Value = "  -1.10%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08431"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.066"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.755"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.92%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.00001307"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.58%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.644.60"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.66%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "94.47"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.46%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06971"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "19.63"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.99%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.856"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  +1.03%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "23.890.43"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.480"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.33%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.977"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +4.61%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "152.81"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.62%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.434"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.57%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "138.75"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.721"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.488"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.824.31"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.025"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +4.77%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.08013"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02942"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.36%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "6.671"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.46%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "10.83"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +5.21%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.2674"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.09097"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.7509"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("E43").Value = "  -1.23%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.23"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.6888"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.432"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.064"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("E48").Value = "  +0.10%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.08268"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "133.97"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.68%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.219"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
